# Add "Chapter 5" (column F) index-term check values for the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = "n/a"
$ws.Range("F4").Value = "n/a"
$ws.Range("F5").Value = "n/a"
$ws.Range("F6").Value = "n/a"
$ws.Range("F7").Value = "check"
$ws.Range("F8").Value = "n/a"
$ws.Range("F9").Value = "n/a"
$ws.Range("F10").Value = "n/a"
$ws.Range("F11").Value = "n/a"
$ws.Range("F12").Value = "n/a"
$ws.Range("F13").Value = "n/a"
$ws.Range("F14").Value = "n/a"
$ws.Range("F15").Value = "n/a"
$ws.Range("F16").Value = "n/a"
$ws.Range("F17").Value = "n/a"
$ws.Range("F18").Value = "n/a"
$ws.Range("F19").Value = "check"
$ws.Range("F20").Value = "check"
$ws.Range("F21").Value = "check"
$ws.Range("F22").Value = "n/a"
$ws.Range("F23").Value = "check"
$ws.Range("F24").Value = "n/a"
$ws.Range("F25").Value = "n/a"
$ws.Range("F26").Value = "n/a"
$ws.Range("F27").Value = "check"
$ws.Range("F28").Value = "n/a"
$ws.Range("F29").Value = "n/a"
$ws.Range("F30").Value = "check"
$ws.Range("F31").Value = "check"
$ws.Range("F32").Value = "n/a"
$ws.Range("F33").Value = "n/a"
$ws.Range("F34").Value = "n/a"
$ws.Range("F35").Value = "n/a"
$ws.Range("F36").Value = "n/a"
$ws.Range("F37").Value = "check"
$ws.Range("F38").Value = "n/a"
$ws.Range("F39").Value = "n/a"
$ws.Range("F40").Value = "n/a"
$ws.Range("F41").Value = "check"
$ws.Range("F42").Value = "check"
$ws.Range("F43").Value = "check"
$ws.Range("F44").Value = "n/a"
$ws.Range("F45").Value = "n/a"
$ws.Range("F46").Value = "check"
$ws.Range("F47").Value = "n/a"
$ws.Range("F48").Value = "n/a"
$ws.Range("F49").Value = "n/a"
$ws.Range("F50").Value = "n/a"
$ws.Range("F51").Value = "n/a"
$ws.Range("F52").Value = "n/a"
$ws.Range("F53").Value = "n/a"
$ws.Range("F54").Value = "check"
$ws.Range("F55").Value = "check"
$ws.Range("F56").Value = "n/a"
$ws.Range("F57").Value = "n/a"
$ws.Range("F58").Value = "n/a"
$ws.Range("F59").Value = "n/a"
$ws.Range("F60").Value = "n/a"
$ws.Range("F61").Value = "n/a"
$ws.Range("F64").Value = "n/a"
$ws.Range("F65").Value = "n/a"
$ws.Range("F66").Value = "n/a"
$ws.Range("F67").Value = "n/a"
$ws.Range("F68").Value = "n/a"
$ws.Range("F69").Value = "n/a"
$ws.Range("F70").Value = "check"
$ws.Range("F71").Value = "check"
$ws.Range("F72").Value = "n/a"
$ws.Range("F73").Value = "n/a"
$ws.Range("F74").Value = "check"
$ws.Range("F75").Value = "check"
$ws.Range("F76").Value = "n/a"
$ws.Range("F77").Value = "n/a"

# Update the visible scroll position / selection to match the saved view.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 58
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F78").Select()

